# constants.xlsx edit script
# Commit: "Fixed fuel mass eq and connected some variables
#          also added solver in correct way, now gets proper convergence but no optimisation (no desvars)"
#
# Semantically this commit:
#   1. Removes the obsolete "T" (thrust) row from the constants table (old row 24),
#      shifting every row below it up by one.
#   2. Renames the description of the (now shifted-up) "S" / wing-area row
#      from "wing area" to "wing area INITIAL " (trailing space kept).
#   3. Corrects the TSFC value in C21 from 5.0458715596330272E-7 to 6.0639E-5.
#   4. Leaves the cursor/selection on C21 (the cell that was actually fixed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the obsolete "T" / thrust row (row 24) -----------------------
# This is a row inside Table1, so Excel automatically shrinks the table
# range/autofilter (A1:H64 -> A1:H63), the sheet dimension, and re-indexes
# the shared strings (dropping now-unused "T"/"thrust"/"kN").
$ws.Rows("24").Delete()

# --- 2. Fix the TSFC value that used to live at C21 --------------------------
$ws.Range("C21").Value = 0.000060639

# --- 3. Rename the wing-area row's description (now row 24 after the delete) -
$ws.Range("B24").Value = "wing area INITIAL "

# --- 4. Keep the _xlnm._FilterDatabase defined name in sync with the table ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$H`$63"
    }
}

# --- 5. Keep the conditional-formatting range in sync with the shrunk table --
$oldCfRange = $ws.Range("E2:H64")
$fcs = $oldCfRange.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("E2:H63"))
}

# --- 6. Leave the selection on the cell that was actually edited -------------
$ws.Range("C21").Select() | Out-Null
